$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 1415
    $ws.Range("F4").Value = 93
    $ws.Range("F6").Value = 9
}
